# "Generate Report for Handoff"
# Updates the localization-status report after a new handoff round:
#   - Overview sheet: refresh "Latest HO Xliff Generate Date" for the
#     files that were just re-handed-off.
#   - zh-cn / de-de sheets: refresh "Latest Handoff Datetime" for those
#     same rows, and mark their "Priority" as "ht" (hot/handoff-triggered).

$wb = $excel.ActiveWorkbook

$rows = @(8, 9, 10, 12, 13, 14)

# --- Overview sheet -------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $overview.Cells.Item($r, 7).Value = "2016-08-30 10:21:22"
}

# --- zh-cn sheet ------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $zhcn.Cells.Item($r, 5).Value = "ht"
    $zhcn.Cells.Item($r, 8).Value = "2016-08-30 10:21:17"
}

# --- de-de sheet ------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $dede.Cells.Item($r, 5).Value = "ht"
    $dede.Cells.Item($r, 8).Value = "2016-08-30 10:21:22"
}
